$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

$ws.Range("A35").Value = "FN0821"
$ws.Range("B35").Value = "Aboriginal"
$ws.Range("C35").Value = "Aug 15-21"
$ws.Range("D35").Value = "Sockeye"
$ws.Range("E35").Value = "Selective Gear"
$ws.Range("F35").Value = "Region 6-LBN"
$ws.Range("G35").Value = 7
$ws.Range("I35").Value = "Sockeye target, selective gear only"

$ws.Range("C35").NumberFormat = $ws.Range("C34").NumberFormat

[void]$ws.Range("C37").Select()
